$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: "administrativo" -> "recursos humanos"; "Compra de Passagem" -> "teste"
$ws.Range("A3").Value = "recursos humanos"
$ws.Range("E3").Value = "teste"

# Row 4: "administrativo" -> "recursos humanos"; new B4 = "r2"; "Sala de Reunião" -> "item1"
$ws.Range("A4").Value = "recursos humanos"
$ws.Range("B4").Value = "r2"
$ws.Range("E4").Value = "item1"

# Rows 5-10 held the old catalog sample rows; content is wiped (formatting/styles kept)
$ws.Range("A5:E10").ClearContents()

# Column A is widened to best-fit the new, longer "recursos humanos" text
$ws.Columns("A").AutoFit() | Out-Null
$ws.Columns("A").ColumnWidth = 15.6

# Selection ends on E5
$ws.Range("E5").Select() | Out-Null
